$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Fill in the new timesheet entries (values only, no formatting
#    yet). Order of entry matters because it controls the order new
#    shared strings get interned in the workbook: A/B/C column block
#    first, then K/L/M column block, then F/G/H column block - this
#    mirrors the order new <si> entries appear in the target file.
# ------------------------------------------------------------------

# --- Column block A:C (rows 22-29) ---
$ws.Cells.Item(22,1).Value = 43052
$ws.Cells.Item(22,2).Value = "Pushing and finished HighScore Unit Test"
$ws.Cells.Item(22,3).Value = 25

$ws.Cells.Item(23,1).Value = 43053
$ws.Cells.Item(23,2).Value = "Team Meeting: Beta Planning"
$ws.Cells.Item(23,3).Value = 40

$ws.Cells.Item(24,1).Value = 43055
$ws.Cells.Item(24,2).Value = "Working on Difficulty, Levels, and Level Timers"
$ws.Cells.Item(24,3).Value = 180

$ws.Cells.Item(25,1).Value = 43056
$ws.Cells.Item(25,2).Value = "Added High Score"
$ws.Cells.Item(25,3).Value = 60

$ws.Cells.Item(26,1).Value = 43057
$ws.Cells.Item(26,2).Value = "Finished High Score, Sound, worked on PowerUps, fixed bugs"
$ws.Cells.Item(26,3).Value = 330

$ws.Cells.Item(27,1).Value = 43057
$ws.Cells.Item(27,2).Value = "Finished Help Screen"
$ws.Cells.Item(27,3).Value = 30

$ws.Cells.Item(28,1).Value = 43057
$ws.Cells.Item(28,2).Value = "Team Meeting"
$ws.Cells.Item(28,3).Value = 60

$ws.Cells.Item(29,1).Value = 43058
$ws.Cells.Item(29,2).Value = "Noah and I talked"
$ws.Cells.Item(29,3).Value = 30

# --- Column block K:M (rows 20-27) ---
$ws.Cells.Item(20,11).Value = 43053
$ws.Cells.Item(20,12).Value = "meeting for beta, implemented load formation"
$ws.Cells.Item(20,13).Value = 40

$ws.Cells.Item(21,11).Value = 43055
$ws.Cells.Item(21,12).Value = "implemented Tracker AI, ad mines with their loads"
$ws.Cells.Item(21,13).Value = 75

$ws.Cells.Item(22,11).Value = 43055
$ws.Cells.Item(22,12).Value = "edited load/save"
$ws.Cells.Item(22,13).Value = 30

$ws.Cells.Item(23,11).Value = 43057
$ws.Cells.Item(23,12).Value = "implemented PowerUp class and its load/save"
$ws.Cells.Item(23,13).Value = 31

$ws.Cells.Item(24,11).Value = 43057
$ws.Cells.Item(24,12).Value = "implemented powerups in player and editted a few things"
$ws.Cells.Item(24,13).Value = 24

$ws.Cells.Item(25,11).Value = 43057
$ws.Cells.Item(25,12).Value = "repaired load/save unit tests"
$ws.Cells.Item(25,13).Value = 71

$ws.Cells.Item(26,11).Value = 43057
$ws.Cells.Item(26,12).Value = "meeting, discussed beta"
$ws.Cells.Item(26,13).Value = 50

$ws.Cells.Item(27,11).Value = 43058
$ws.Cells.Item(27,12).Value = "made Beta video"
$ws.Cells.Item(27,13).Value = 35

# --- Column block F:H (rows 24-28) ---
$ws.Cells.Item(24,6).Value = 43053
$ws.Cells.Item(24,7).Value = "Meeting about beta"
$ws.Cells.Item(24,8).Value = 40

$ws.Cells.Item(25,6).Value = 43055
$ws.Cells.Item(25,7).Value = "Talked with Robert about level design"
$ws.Cells.Item(25,8).Value = 5

$ws.Cells.Item(26,6).Value = 43057
$ws.Cells.Item(26,7).Value = "Added death"
$ws.Cells.Item(26,8).Value = 40

$ws.Cells.Item(27,6).Value = 43057
$ws.Cells.Item(27,7).Value = "team meeting"
$ws.Cells.Item(27,8).Value = 55

$ws.Cells.Item(28,6).Value = 43058
$ws.Cells.Item(28,7).Value = "fixed powerup spawn on load"
$ws.Cells.Item(28,8).Value = 15

# ------------------------------------------------------------------
# 2) Now that all values are in place, extend the alternating-row
#    formatting (fills, fonts, number formats keyed off style indices
#    1-9) down into the new rows by pasting *formats only* from the
#    existing odd/even template rows. PasteSpecial(xlPasteFormats)
#    does not touch cell values, so this is safe to run after filling
#    in the data above.
# ------------------------------------------------------------------
$xlPasteFormats = -4122

# Columns A:C -> replicate rows 20:21 pattern down through rows 22:29
$ws.Range("A20:C21").Copy()
$ws.Range("A22:C23").PasteSpecial($xlPasteFormats)
$ws.Range("A24:C25").PasteSpecial($xlPasteFormats)
$ws.Range("A26:C27").PasteSpecial($xlPasteFormats)
$ws.Range("A28:C29").PasteSpecial($xlPasteFormats)

# Columns F:H -> replicate rows 20:21 pattern down through rows 24:28
$ws.Range("F20:H21").Copy()
$ws.Range("F24:H25").PasteSpecial($xlPasteFormats)
$ws.Range("F26:H27").PasteSpecial($xlPasteFormats)
$ws.Range("F20:H20").Copy()
$ws.Range("F28:H28").PasteSpecial($xlPasteFormats)

# Columns K:M -> replicate rows 18:19 pattern down through rows 20:27
$ws.Range("K18:M19").Copy()
$ws.Range("K20:M21").PasteSpecial($xlPasteFormats)
$ws.Range("K22:M23").PasteSpecial($xlPasteFormats)
$ws.Range("K24:M25").PasteSpecial($xlPasteFormats)
$ws.Range("K26:M27").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3) Recalculate so cached formula results (SUM / derived cells) match.
# ------------------------------------------------------------------
$excel.CalculateFull()

# ------------------------------------------------------------------
# 4) Update the view/selection to match the edited workbook state.
# ------------------------------------------------------------------
$ws.Range("I11").Select()
